# Update the "EmailSent"/"Response" columns (X/Y) into a single
# "DataOrigin" column (X) that records where each study's extracted
# data came from, and drop the now-unused "Response" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header and rewrite the data values in column X.
$ws.Range("X1").Value = "DataOrigin"

$ws.Range("X2").Value  = "Figure 2"
$ws.Range("X3").Value  = "Table 2"
$ws.Range("X4").Value  = "Figure 3e"
$ws.Range("X5").Value  = "Table 3"
$ws.Range("X6").Value  = "Figure 2,3,4"
$ws.Range("X7").Value  = "Figure 2b"
$ws.Range("X8").Value  = "Figure 1"
$ws.Range("X9").Value  = "Figure 2"
$ws.Range("X10").Value = "Figure 1,3a,4"
$ws.Range("X11").Value = "Figure 4, Table 1"

# The old "Response" column (Y) was always empty below the header, so
# simply delete the whole column, shrinking the used range back to A:X.
$ws.Range("Y1:Y11").EntireColumn.Delete()

# Refresh the view state to match: frozen pane now starts at Q2 and the
# active selection sits on X3.
$ws.Range("Q2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("X3").Select()
